$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.248.26"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "2.616.25"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'523.72"
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("D6").Value = "'144.57"
$ws.Range("E6").Value = "  +1.45%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "2.614.50"
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").Value = "'6.65"
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("E11").Value = "  -0.89%  "
$ws.Range("D12").Value = "'0.335"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("E13").Value = "  -0.63%  "
$ws.Range("D14").Value = "3.072.10"
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("D15").Value = "58.204.75"
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("E16").Value = "  -1.80%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.0000134"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.625.04"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").Value = "'340.16"
$ws.Range("E19").Value = "  +1.70%  "
$ws.Range("D20").Value = "'4.37"
$ws.Range("E20").Value = "  -0.62%  "
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("D22").Value = "'6.39"
$ws.Range("E22").Value = "  +2.26%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").Value = "'65.51"
$ws.Range("E24").Value = "  +2.42%  "
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("E26").Value = "  -2.54%  "
$ws.Range("D27").Value = "2.716.45"
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("D28").Value = "'0.997"
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("D29").Value = "'7.04"
$ws.Range("E29").Value = "  -0.76%  "
$ws.Range("D30").Value = "0.0₃0751"
$ws.Range("E30").Value = "  -4.35%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("E32").Value = "  -5.11%  "
$ws.Range("E33").Value = "  +0.82%  "
$ws.Range("E34").Value = "  +1.07%  "
$ws.Range("D35").Value = "'149.82"
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("E36").Value = "  -0.59%  "
$ws.Range("E37").Value = "  -2.00%  "
$ws.Range("D38").Value = "'0.866"
$ws.Range("E38").Value = "  -3.66%  "
$ws.Range("D39").Value = "'0.851"
$ws.Range("E39").Value = "  +1.41%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "'36.12"
$ws.Range("E40").Value = "  -1.33%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.46"
$ws.Range("E41").Value = "  +2.17%  "
$ws.Range("E42").Value = "  -1.19%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Value = "'273.65"
$ws.Range("E44").Value = "  +2.11%  "
$ws.Range("D45").Value = "'0.597"
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("E46").Value = "  -0.50%  "
$ws.Range("E47").Value = "  +0.39%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'18.86"
$ws.Range("E48").Value = "  -1.34%  "
$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").Value = "'0.0525"
$ws.Range("E49").Value = "  -1.30%  "
$ws.Range("D50").Value = "'19.18"
$ws.Range("E50").Value = "  +5.26%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "'4.67"
$ws.Range("E51").Value = "  +1.00%  "
